# "modell was split into import and analysis"
# Adds a "bus" and "marginal_cost" row to the technical-data tables on both
# sheets, and bumps the installed-capacity (D7) figures used for the model
# import.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("GuD Kraftwerk")
$ws2 = $wb.Worksheets.Item("Gasturbinenkraftwerk")

# --- GuD Kraftwerk ------------------------------------------------------
$ws1.Range("D7").Value = 3500

$ws1.Range("C13").Value = "bus"
$ws1.Range("D13").Value = 2

$ws1.Range("C14").Value = "marginal_cost"
$ws1.Range("D14").Value = 30

# --- Gasturbinenkraftwerk ------------------------------------------------
$ws2.Range("D7").Value = 3000

$ws2.Range("C13").Value = "bus"
$ws2.Range("D13").Value = 1

$ws2.Range("C14").Value = "marginal_cost"
$ws2.Range("D14").Value = 50

# --- Selection: land on C15 on both sheets, leaving Gasturbinenkraftwerk
# (the originally active tab) active when the script finishes. -----------
$ws1.Range("C15").Select()
$ws2.Activate()
$ws2.Range("C15").Select()
